$d = $word.ActiveDocument

# The title/application-note paragraph ("Application Note: Assembling freETarget Competition")
# is the second paragraph in the document.
$p = $d.Paragraphs.Item(2)
$r = $p.Range

# Locate "Competition" within that paragraph so we can append " and Club" right after it,
# mirroring a user placing their cursor at the end of "Competition" and typing more text.
$find = $r.Duplicate
$find.Find.ClearFormatting()
$found = $find.Find.Execute("Competition", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $ins = $find.Duplicate
    $ins.Collapse(0)  # wdCollapseEnd
    $ins.InsertAfter(" and Club")
}

# Re-fetch the paragraph's full range (now including the newly-typed text) and shrink the
# whole paragraph from 20pt down to 18pt, matching the size reduction applied across every
# run in the paragraph (including the paragraph mark and the anchored drawing's run).
$p2 = $d.Paragraphs.Item(2)
$full = $p2.Range
$full.Font.Size = 18
$full.Font.SizeBi = 18
